# Rename the existing "Users" sheet to "ValidUsers"
$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item("Users")
$wsValid.Name = "ValidUsers"

# Add a new worksheet "InvalidUsers" right after "ValidUsers" with an invalid login test case
$wsInvalid = $wb.Worksheets.Add($null, $wsValid)
$wsInvalid.Name = "InvalidUsers"

$wsInvalid.Range("A1").Value = "userName"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "admin"
$wsInvalid.Range("B2").Value = "trainee"
$wsInvalid.Range("A3").Value = "trainee"
$wsInvalid.Range("B3").Value = "manager"

# The rows that carry the new invalid-login data keep an explicit row height,
# matching how Excel marks rows touched interactively while entering data.
$wsInvalid.Rows.Item(2).RowHeight = 15.75
$wsInvalid.Rows.Item(3).RowHeight = 15.75

$null = $wsInvalid.Range("A2").Select()

# Re-select ValidUsers and set the active cell to A3 (matches target file)
$null = $wsValid.Activate()
$null = $wsValid.Range("A3").Select()
